# ACDynamicStabilityManager balance update:
# The Xcg "ESTIMATION METHOD COMPARISON" table rows for SFORZA / TORENBEEK_1982
# were reshuffled across several sheets.
#
# Worksheets collection order (1-based), matching workbook.xml sheet order:
#   1  GLOBAL RESULTS
#   2  FUSELAGE
#   3  WING
#   4  FUEL TANK
#   5  HORIZONTAL TAIL
#   6  VERTICAL TAIL
#   7  CANARD
#   8  NACELLES
#   9  POWER PLANT
#   10 LANDING GEARS
#   11 SYSTEMS

$wb = $excel.ActiveWorkbook

# --- FUSELAGE: the SFORZA and TORENBEEK_1982 rows (label + value) swap places ---
$wsFuselage = $wb.Worksheets.Item(2)
$wsFuselage.Range("A23").Value() = "TORENBEEK_1982"
$wsFuselage.Range("C23").Value() = 17.878799999999995
$wsFuselage.Range("A24").Value() = "SFORZA"
$wsFuselage.Range("C24").Value() = 15.527255597045638

# --- WING: same swap, occurring twice (Xcg block rows 23/24 and Ycg block rows 27/28) ---
$wsWing = $wb.Worksheets.Item(3)
$wsWing.Range("A23").Value() = "TORENBEEK_1982"
$wsWing.Range("C23").Value() = 1.8563139127409203
$wsWing.Range("A24").Value() = "SFORZA"
$wsWing.Range("C24").Value() = 2.2926530320804064
$wsWing.Range("A27").Value() = "TORENBEEK_1982"
$wsWing.Range("C27").Value() = 6.868233768733532
$wsWing.Range("A28").Value() = "SFORZA"
$wsWing.Range("C28").Value() = 7.957351173907407

# --- HORIZONTAL TAIL: rows already use the TORENBEEK_1982 method, relabel to keep in sync ---
$wsHTail = $wb.Worksheets.Item(5)
$wsHTail.Range("A23").Value() = "TORENBEEK_1982"
$wsHTail.Range("A26").Value() = "TORENBEEK_1982"

# --- VERTICAL TAIL: same as horizontal tail ---
$wsVTail = $wb.Worksheets.Item(6)
$wsVTail.Range("A23").Value() = "TORENBEEK_1982"
$wsVTail.Range("A26").Value() = "TORENBEEK_1982"

# --- CANARD: same as horizontal tail ---
$wsCanard = $wb.Worksheets.Item(7)
$wsCanard.Range("A23").Value() = "TORENBEEK_1982"
$wsCanard.Range("A26").Value() = "TORENBEEK_1982"

# --- POWER PLANT: row uses the SFORZA method ---
$wsPowerPlant = $wb.Worksheets.Item(9)
$wsPowerPlant.Range("A23").Value() = "SFORZA"

# --- LANDING GEARS: rows use the SFORZA method ---
$wsLandingGears = $wb.Worksheets.Item(10)
$wsLandingGears.Range("A23").Value() = "SFORZA"
$wsLandingGears.Range("A26").Value() = "SFORZA"
